# Fixes for race condition
# Rework "Prev" row label -> "P", renumber shared strings implicitly,
# and add a second (LHS/Direction) verification table in rows 8-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing "Prev" label cell to "P" ---
$ws.Range("A2").Value = "P"

# --- New column widths / formatting for columns A-D ---
$ws.Columns("A").ColumnWidth = 3.11
$ws.Columns("B").ColumnWidth = 2.66
$ws.Columns("C").ColumnWidth = 2.66
$ws.Columns("D").ColumnWidth = 4.11
$ws.Columns("E").ColumnWidth = 16.44

# Bold font for column A labels (A2:A5)
$ws.Range("A2:A5").Font.Bold = $true
$ws.Range("A2:A5").HorizontalAlignment = -4152  # xlRight

# --- New table header row (row 8) ---
$ws.Range("B8").Value = "I"
$ws.Range("C8").Value = "Q"
$ws.Range("D8").Value = "LHS"
$ws.Range("E8").Value = "Direction"

$headerRange = $ws.Range("A8:E8")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$ws.Range("A8:D8").HorizontalAlignment = -4152  # xlRight
$ws.Range("E8").HorizontalAlignment = -4131     # xlLeft

# --- New rows 9-12: labels + formulas ---
$ws.Range("A9").Value = "P"
$ws.Range("A10").Value = "C0"
$ws.Range("A11").Value = "C1"
$ws.Range("A12").Value = "C2"

$ws.Range("A9:A12").Font.Bold = $true

$ws.Range("B9").Formula = "=C2"
$ws.Range("C9").Formula = "=B2"

$ws.Range("B10").Formula = "=C3"
$ws.Range("C10").Formula = "=B3"
$ws.Range("D10").Formula = "=F3"
$ws.Range("E10").Value = "counterclockwise"

$ws.Range("B11").Formula = "=C4"
$ws.Range("C11").Formula = "=B4"
$ws.Range("D11").Formula = "=F4"
$ws.Range("E11").Value = "clockwise"

$ws.Range("B12").Formula = "=C5"
$ws.Range("C12").Formula = "=B5"
$ws.Range("D12").Formula = "=F5"
$ws.Range("E12").Value = "counterclockwise"

# Borders around the whole new table (rows 8-12, cols A-E)
$tableRange = $ws.Range("A8:E12")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

$ws.Range("E9:E12").NumberFormat = $ws.Range("E2").NumberFormat

# --- Selection ---
$ws.Range("A8:E12").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1  # xlPortrait
